# Generate Report for Handback
#
# The "97814758-b34a-46ff-8abf-bdbd317417c5.md" file has completed its
# handback cycle: its Status flips from "Ready for handoff" to
# "Handed back: in sync with en-US" on the Overview sheet as well as on
# each per-locale sheet, and the per-locale "Latest Handback DateTime"
# stamps are refreshed to reflect the new handback.

$wb = $excel.ActiveWorkbook

$statusDone = "Handed back: in sync with en-US"

# --- Overview sheet -------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $statusDone
$overview.Range("C3").Value = $statusDone

# --- zh-cn sheet ------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $statusDone
$zhcn.Range("H3").Value = "2016-03-18 05:47:24"

# --- de-de sheet ------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $statusDone
$dede.Range("H3").Value = "2016-03-18 05:47:29"
